$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.309.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.85"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07883"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.24"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08397"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.877.71"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.245"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7179"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.21"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.209"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008356"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.308.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.82"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.122.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.796"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1594"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.76"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.056"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.54"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.507"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.424"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.16%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.209"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.14%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.19%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7512"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.178"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.698"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.292.59"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.77%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.731"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.573"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "110.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8956"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "73.29"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000130"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +8.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.012.14"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.803"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5203"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.454"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4360"
